# "Add code for graph creation and data processing"
#
# Data edits:
#   - node sheet: the totals row (row 7, holding SUM() formulas) is removed.
#   - node sheet: the "D" row's s_clean/s_dirty figures (C5/D5) are swapped
#     (160 and 10 traded places).
#   - edge sheet: the "cost" column header is renamed to "price".
#   - the "edge" sheet becomes the active/front tab, with its selection left
#     on the freshly edited header cell; the "node" sheet's selection is left
#     on the now-cleared former totals row.

$wb = $excel.ActiveWorkbook

$wsNode = $wb.Worksheets.Item("node")
$wsEdge = $wb.Worksheets.Item("edge")

# --- node: drop the SUM() totals row entirely ---
$wsNode.Range("B7:E7").ClearContents()

# --- node: swap the s_clean / s_dirty values for row 5 ---
$wsNode.Range("C5").Value = 160
$wsNode.Range("D5").Value = 10

# leave node's selection where the totals row used to be
$wsNode.Range("B7:E7").Select()

# --- edge: rename the "cost" header to "price" ---
$wsEdge.Range("C1").Value = "price"

# edge becomes the active sheet, selection resting on the renamed header
$wsEdge.Activate()
$wsEdge.Range("C1").Select()
